$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 129
$ws.Range("I6").Value = 129
$ws.Range("K6").Value = 387
$ws.Range("M6").Value = -275
$ws.Range("H98").Value = 1245.0625
$ws.Range("I98").Value = 1287.2142
$ws.Range("J98").Value = 950
$ws.Range("K98").Value = 1287.2142
$ws.Range("L98").Value = 950
$ws.Range("M98").Value = 210.7858000000001
$ws.Range("N98").Value = -3946
$ws.Range("H100").Value = 2993.4167
$ws.Range("J100").Value = 3367.889
$ws.Range("L100").Value = 3367.889
$ws.Range("N100").Value = -4449.889
$ws.Range("H112").Value = 3978.4375
$ws.Range("J112").Value = 3978.4375
$ws.Range("L112").Value = 11935.3125
$ws.Range("N112").Value = -14151.3125
$ws.Range("H113").Value = 202651.7
$ws.Range("I113").Value = 2905
$ws.Range("J113").Value = 224845.78
$ws.Range("K113").Value = 2905
$ws.Range("L113").Value = 224845.78
$ws.Range("M113").Value = 349
$ws.Range("N113").Value = -231353.78
$ws.Range("H122").Value = 1245.0625
$ws.Range("I122").Value = 1287.2142
$ws.Range("J122").Value = 950
$ws.Range("K122").Value = 3861.6426
$ws.Range("L122").Value = 2850
$ws.Range("M122").Value = -1411.6426
$ws.Range("N122").Value = -7750
$ws.Range("H127").Value = 1140.3334
$ws.Range("I127").Value = 1260.75
$ws.Range("K127").Value = 3782.25
$ws.Range("M127").Value = 1177.75
$ws.Range("H131").Value = 2424.75
$ws.Range("I131").Value = 2424.75
$ws.Range("K131").Value = 7274.25
$ws.Range("M131").Value = -2234.25
$ws.Range("H132").Value = 1202.9474
$ws.Range("I132").Value = 991.58826
$ws.Range("K132").Value = 2974.76478
$ws.Range("M132").Value = -444.76478
$ws.Range("H135").Value = 2513.6785
$ws.Range("I135").Value = 1747.4375
$ws.Range("J135").Value = 3535.3333
$ws.Range("K135").Value = 15726.9375
$ws.Range("L135").Value = 31817.9997
$ws.Range("M135").Value = -13191.9375
$ws.Range("N135").Value = -36887.9997
$ws.Range("H137").Value = 6953005.5
$ws.Range("I137").Value = 13162389
$ws.Range("K137").Value = 39487167
$ws.Range("M137").Value = -39484617
$ws.Range("H138").Value = 3249.775
$ws.Range("J138").Value = 3623.1072
$ws.Range("L138").Value = 10869.3216
$ws.Range("N138").Value = -21149.3216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 23736.5
$ws.Range("J45").Value = 1648.6666
$ws.Range("L45").Value = 1648.6666
$ws.Range("N45").Value = -2402.6666
$ws.Range("H61").Value = 4239.7915
$ws.Range("I61").Value = 3631.6667
$ws.Range("J61").Value = 8496.666999999999
$ws.Range("K61").Value = 3631.6667
$ws.Range("L61").Value = 8496.666999999999
$ws.Range("M61").Value = -3419.6667
$ws.Range("N61").Value = -8920.666999999999
$ws.Range("H63").Value = 8082.4053
$ws.Range("J63").Value = 11005.792
$ws.Range("L63").Value = 11005.792
$ws.Range("N63").Value = -12377.792
$ws.Range("H66").Value = 8082.4053
$ws.Range("J66").Value = 11005.792
$ws.Range("L66").Value = 55028.96
$ws.Range("N66").Value = -61892.96
$ws.Range("H97").Value = 1899980.8
$ws.Range("I97").Value = 2473163.5
$ws.Range("J97").Value = 180432.2
$ws.Range("K97").Value = 2473163.5
$ws.Range("L97").Value = 180432.2
$ws.Range("M97").Value = -2472667.5
$ws.Range("N97").Value = -181424.2
$ws.Range("H110").Value = 12500803
$ws.Range("I110").Value = 12500803
$ws.Range("K110").Value = 12500803
$ws.Range("M110").Value = -12498758
$ws.Range("H132").Value = 4726.278
$ws.Range("I132").Value = 3290.9285
$ws.Range("K132").Value = 9872.7855
$ws.Range("M132").Value = -7342.7855
$ws.Range("H136").Value = 4239.7915
$ws.Range("I136").Value = 3631.6667
$ws.Range("J136").Value = 8496.666999999999
$ws.Range("K136").Value = 10895.0001
$ws.Range("L136").Value = 25490.001
$ws.Range("M136").Value = -8345.000100000001
$ws.Range("N136").Value = -30590.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 184817.64
$ws.Range("I86").Value = 2999.1428
$ws.Range("K86").Value = 2999.1428
$ws.Range("M86").Value = -1876.1428
$ws.Range("H89").Value = 184817.64
$ws.Range("I89").Value = 2999.1428
$ws.Range("K89").Value = 14995.714
$ws.Range("M89").Value = -9379.714
$ws.Range("H94").Value = 789.55554
$ws.Range("I94").Value = 708.8461
$ws.Range("K94").Value = 708.8461
$ws.Range("M94").Value = -257.8461

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 33334334
$ws.Range("J31").Value = 4794.0625
$ws.Range("K31").Value = 33334334
$ws.Range("L31").Value = 4794.0625
$ws.Range("M31").Value = -33334039
$ws.Range("N31").Value = -5384.0625
$ws.Range("I34").Value = 33334334
$ws.Range("J34").Value = 4794.0625
$ws.Range("K34").Value = 33334334
$ws.Range("L34").Value = 4794.0625
$ws.Range("M34").Value = -33334132
$ws.Range("N34").Value = -5198.0625
$ws.Range("H52").Value = 79000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 79000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 79000
$ws.Range("N52").Value = -79588
$ws.Range("M52").ClearContents()
$ws.Range("H58").Value = 3058.52
$ws.Range("I58").Value = 1809.2778
$ws.Range("J58").Value = 6270.857
$ws.Range("K58").Value = 1809.2778
$ws.Range("L58").Value = 6270.857
$ws.Range("M58").Value = -1606.2778
$ws.Range("N58").Value = -6676.857
$ws.Range("H125").Value = 85348.625
$ws.Range("J125").Value = 85348.625
$ws.Range("L125").Value = 85348.625
$ws.Range("N125").Value = -90268.625
$ws.Range("H132").Value = 144496.77
$ws.Range("I132").Value = 100774.11
$ws.Range("K132").Value = 302322.33
$ws.Range("M132").Value = -299792.33
$ws.Range("H133").Value = 61498
$ws.Range("J133").Value = 61498
$ws.Range("L133").Value = 61498
$ws.Range("N133").Value = -66558
$ws.Range("H134").Value = 12949.2
$ws.Range("I134").Value = 14499.5
$ws.Range("J134").Value = 11915.667
$ws.Range("K134").Value = 43498.5
$ws.Range("L134").Value = 35747.001
$ws.Range("M134").Value = -40963.5
$ws.Range("N134").Value = -40817.001
$ws.Range("H136").Value = 3058.52
$ws.Range("I136").Value = 1809.2778
$ws.Range("J136").Value = 6270.857
$ws.Range("K136").Value = 5427.8334
$ws.Range("L136").Value = 18812.571
$ws.Range("M136").Value = -2877.8334
$ws.Range("N136").Value = -23912.571
$ws.Range("H137").Value = 114666.336
$ws.Range("J137").Value = 114666.336
$ws.Range("L137").Value = 114666.336
$ws.Range("N137").Value = -124866.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 2395.8333
$ws.Range("I93").Value = 750
$ws.Range("J93").Value = 2944.4443
$ws.Range("K93").Value = 2250
$ws.Range("L93").Value = 8833.332900000001
$ws.Range("M93").Value = -378
$ws.Range("N93").Value = -12577.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 49999
$ws.Range("J74").Value = 49999
$ws.Range("L74").Value = 49999
$ws.Range("N74").Value = -51871
$ws.Range("H77").Value = 49999
$ws.Range("J77").Value = 49999
$ws.Range("L77").Value = 149997
$ws.Range("N77").Value = -159357

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H94").Value = 55895
$ws.Range("J94").Value = 55895
$ws.Range("L94").Value = 55895
$ws.Range("N94").Value = -57247
$ws.Range("H136").Value = 4669.1035
$ws.Range("I136").Value = 3054.9524
$ws.Range("K136").Value = 9164.8572
$ws.Range("M136").Value = -6614.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5743.091
$ws.Range("I132").Value = 3288.125
$ws.Range("K132").Value = 9864.375
$ws.Range("M132").Value = -7334.375
$ws.Range("H136").Value = 2237.7144
$ws.Range("I136").Value = 1294.4615
$ws.Range("K136").Value = 3883.3845
$ws.Range("M136").Value = -1333.3845
